# CompressionResults.xlsx - "Finished!" commit
# - Fixes LZWmod (no reset) / LZWmod (with reset) raw byte counts for several files
# - Flips the compression-ratio formulas to Original/Compressed (using table
#   structured references instead of raw cell refs) across columns E, G, I, K
# - Re-formats the ratio columns' numbers (0.000 -> 0.00, #,##0.000 -> #,##0.00)
# - Leaves the selection parked on K20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct the handful of raw byte counts that were wrong before (diff shows
#    these as plain <v> changes with no formula attached).
# ---------------------------------------------------------------------------
$ws.Range("F4").Value  = 1792781
$ws.Range("H4").Value  = 1177887

$ws.Range("F9").Value  = 156409
$ws.Range("H9").Value  = 152231

$ws.Range("H10").Value = 171170

$ws.Range("F12").Value = 501777
$ws.Range("H12").Value = 527598

$ws.Range("F14").Value = 597847
$ws.Range("H14").Value = 590545

# ---------------------------------------------------------------------------
# 2. Rewrite the compression-ratio formulas (columns E, G, I, K, rows 3-16) so
#    they read Original/Compressed via Table1 structured references instead
#    of the old Compressed/Original raw-cell-ref formulas.
# ---------------------------------------------------------------------------
$ratioCols = @{
    "E" = "LZW"
    "G" = "LZWmod (no reset)"
    "I" = "LZWmod (with reset)"
    "K" = "Unix Compress"
}

foreach ($col in $ratioCols.Keys) {
    $colName = $ratioCols[$col]
    $formula = "=Table1[[#This Row],[Original]]/Table1[[#This Row],[$colName]]"
    for ($r = 3; $r -le 16; $r++) {
        $ws.Range("$col$r").Formula = $formula
    }
}

# ---------------------------------------------------------------------------
# 3. Number formats: the "0.000" / "#,##0.000" custom formats become the
#    built-in "0.00" / "#,##0.00" formats.
# ---------------------------------------------------------------------------
$ws.Range("E3:E16").NumberFormat = "0.00"
$ws.Range("G3:G16").NumberFormat = "0.00"
$ws.Range("K3:K16").NumberFormat = "0.00"
$ws.Range("I3:I16").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------------
# 4. Selection moves off the table onto K20.
# ---------------------------------------------------------------------------
$ws.Range("K20").Select()

Write-Output "edit applied"
